# Weekly data refresh: a new day's price record (2021-11-18, serial 44518)
# is inserted at the top of the Albahaca / Terminal La Palmera de La Serena
# block (row 39), pushing every existing row in the block down by one.
# The row that falls off the bottom of the block (the old row 57, serial
# 44432 / 2021-08-24) re-appears as the new last row (58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 39; this shifts rows 39:57 down to 40:58
# and carries the existing formatting (incl. the date number format on
# column D) along with it.
$ws.Rows.Item(39).Insert()

# Populate the newly-inserted row 39 with the new week's record. The
# "boilerplate" columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical for every
# row in this block, so copy them straight from the row below (40), which
# now holds what used to be row 39.
$ws.Cells.Item(39, 1).Value2  = $ws.Cells.Item(40, 1).Value2   # A: Mercado ID
$ws.Cells.Item(39, 2).Value2  = $ws.Cells.Item(40, 2).Value2   # B: Mercado
$ws.Cells.Item(39, 3).Value2  = $ws.Cells.Item(40, 3).Value2   # C: Region
$ws.Cells.Item(39, 4).Value2  = 44518                          # D: Fecha
$ws.Cells.Item(39, 5).Value2  = $ws.Cells.Item(40, 5).Value2   # E: Codreg
$ws.Cells.Item(39, 6).Value2  = $ws.Cells.Item(40, 6).Value2   # F: Categoria ID
$ws.Cells.Item(39, 7).Value2  = $ws.Cells.Item(40, 7).Value2   # G: Categoria
$ws.Cells.Item(39, 8).Value2  = $ws.Cells.Item(40, 8).Value2   # H: Variedad
$ws.Cells.Item(39, 9).Value2  = $ws.Cells.Item(40, 9).Value2   # I: Calidad
$ws.Cells.Item(39, 10).Value2 = 760                             # J: Volumen
$ws.Cells.Item(39, 11).Value2 = 3000                            # K: Precio minimo
$ws.Cells.Item(39, 12).Value2 = 4000                            # L: Precio maximo
$ws.Cells.Item(39, 13).Value2 = 3500                            # M: Precio promedio ponderado
$ws.Cells.Item(39, 14).Value2 = $ws.Cells.Item(40, 14).Value2  # N: Unidad de comercializacion
$ws.Cells.Item(39, 15).Value2 = $ws.Cells.Item(40, 15).Value2  # O: Origen
$ws.Cells.Item(39, 16).Value2 = 3500                            # P: Precio $/Kg
$ws.Cells.Item(39, 17).Value2 = $ws.Cells.Item(40, 17).Value2  # Q: Kg o Unidades
$ws.Cells.Item(39, 18).Value2 = $ws.Cells.Item(40, 18).Value2  # R: Clasificacion
